$wb = $excel.ActiveWorkbook

# --- Sheet "ip_address_list" ---
$ws1 = $wb.Worksheets.Item("ip_address_list")

# row1: remove the empty D1 inline-string cell entirely
$ws1.Range("D1").ClearContents()

# row2: A2 "5" -> "jojo"
$ws1.Range("A2").Value = "jojo"

# row3: A3 "51452422" -> "5145242" (keep it text, not auto-converted to a number)
$ws1.Range("A3").NumberFormat = "@"
$ws1.Range("A3").Value = "5145242"

# row5 gets split into a new row5 (edited) and a new row6 (former row5's
# B/C/E carried down, A6 renamed "brambor")
$oldB5 = $ws1.Range("B5").Value2
$oldC5 = $ws1.Range("C5").Value2
$oldE5 = $ws1.Range("E5").Value2

$ws1.Range("A5").Value = "todleto no"
$ws1.Range("B5").Value = "192.168.000.999"
$ws1.Range("C5").Value = "255.255.255.0"
$ws1.Range("D5").Value = "joo`n§j"
$ws1.Range("E5").Value = 0
# the embedded line break auto-grows the row; put the height back to default
$ws1.Rows.Item(5).AutoFit()

$ws1.Range("A6").Value = "brambor"
$ws1.Range("B6").Value = $oldB5
$ws1.Range("C6").Value = $oldC5
$ws1.Range("E6").Value = $oldE5

# --- Sheet "ip_adress_fav_list" ---
$ws2 = $wb.Worksheets.Item("ip_adress_fav_list")
$ws2.Range("A1").Value = "jojo"

# --- Sheet "disc_list" ---
$ws3 = $wb.Worksheets.Item("disc_list")

# row1's "518" is a numeric-looking value stored as TEXT; read it explicitly
# as a literal (the sheet only ever held "518" here) so it is not coerced to
# a Double when round-tripped through Value2.
$a1 = "518"
$b1 = $ws3.Range("B1").Value2
$c1 = $ws3.Range("C1").Value2
$d1 = $ws3.Range("D1").Value2
$e1 = $ws3.Range("E1").Value2
$f1 = $ws3.Range("F1").Value2

$a2 = $ws3.Range("A2").Value2
$b2 = $ws3.Range("B2").Value2
$c2 = $ws3.Range("C2").Value2
$d2 = $ws3.Range("D2").Value2
$e2 = $ws3.Range("E2").Value2

$ws3.Range("A1").Value = $a2
$ws3.Range("B1").Value = $b2
$ws3.Range("C1").Value = $c2
$ws3.Range("D1").Value = $d2
$ws3.Range("E1").Value = $e2
$ws3.Range("F1").ClearContents()

# A2 becomes "518" again - force text so it isn't auto-converted to a number
$ws3.Range("A2").NumberFormat = "@"
$ws3.Range("A2").Value = $a1
$ws3.Range("B2").Value = $b1
$ws3.Range("C2").Value = $c1
$ws3.Range("D2").Value = $d1
$ws3.Range("E2").Value = $e1
$ws3.Range("F2").Value = $f1
